$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Generate Report for Handback
#
# Updates the localization-status report with the outcome of the
# handback run: the "Ready for handoff" status becomes "Handed back:
# in sync with en-US" everywhere it is shown, and the per-locale
# "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns on the zh-cn and de-de detail sheets are filled in
# with the generated handback artifacts and timestamps.
# ---------------------------------------------------------------------

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$mdFileName  = "7eafc53e-b21e-443d-b489-0c20752e608d.md"
$mdUrl       = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/295b8795d022a9fd64d1cfe9af51fe3c51161df9/e2e/7eafc53e-b21e-443d-b489-0c20752e608d.md"
$statusText  = "Handed back: in sync with en-US"

# 1) Flip the handoff/handback status text everywhere it appears.
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$zhcn.Range("C2").Value     = $statusText
$dede.Range("C2").Value     = $statusText

# 2) zh-cn detail row: target file, handback file and handback datetime.
#    Adding the hyperlink sets the cell's text + the "HyperLink" look in
#    one call.
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl, "", "", $mdFileName) | Out-Null
$zhcn.Range("J2").Value = "7eafc53e-b21e-443d-b489-0c20752e608d.4bf2d6cb022035ed391ddf1fc2651833277492cc.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-21 13:03:00"

# 3) de-de detail row: target file, handback file and handback datetime.
$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl, "", "", $mdFileName) | Out-Null
$dede.Range("J2").Value = "7eafc53e-b21e-443d-b489-0c20752e608d.4bf2d6cb022035ed391ddf1fc2651833277492cc.de-de.xlf"
$dede.Range("K2").Value = "2016-08-21 13:03:11"

# 4) Widen the "Status" column on Overview and the target/handback
#    columns on the per-locale sheets to fit the longer text now shown.
$overview.Range("E1").ColumnWidth = 29.14
$overview.Range("F1").ColumnWidth = 29.14

$zhcn.Range("C1").ColumnWidth = 29.14
$zhcn.Range("I1").ColumnWidth = 39.15
$zhcn.Range("J1").ColumnWidth = 39.15

$dede.Range("C1").ColumnWidth = 29.14
$dede.Range("I1").ColumnWidth = 39.15
$dede.Range("J1").ColumnWidth = 39.15
